$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the claim numbers (NroSiniestro column, F2/F3). Use a leading
# apostrophe so Excel keeps storing them as text (preserving the existing
# text-quoted cell style) instead of reinterpreting them as numbers.
$ws.Range("F2").Value = "'1120194100404"
$ws.Range("F3").Value = "'1120170200928  "

# Update the active selection to match the new cursor position
$ws.Range("F4").Select()
